$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest Pick 3 draw result as a new row (row 60).
# A leading apostrophe forces Excel to treat date-/number-looking values
# ("2025-11-15", "251115", the ISO timestamp) as plain text, matching the
# existing rows in the sheet which are all stored as text.
$ws.Range("A60").Value = "'2025-11-15"
$ws.Range("B60").Value = "Pick 3"
$ws.Range("C60").Value = "'251115"
$ws.Range("D60").Value = "7-0-8"
$ws.Range("E60").Value = "'2025-11-15T21:35:31.380+04:00"

# Keep the "number stored as text" error-checking suppression in sync with
# the newly extended data range (best effort; harmless if unsupported).
try {
    $ws.Range("A60:E60").Errors.Item(9).Ignore = $true
} catch {
}
